$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column B (21.33203125 -> 36.109375)
$ws.Columns.Item(2).ColumnWidth = 36.109375

# C13: "CADA 5 EUROS"
$ws.Range("C13").Value = "CADA 5 EUROS"

# Row 16 headers (typed first, so they land earlier in the shared-string table)
$ws.Range("A16").Value = "IMPORTE CUENTA"
$ws.Range("B16").Value = "PUNTOS"
$ws.Range("C16").Value = "TOTAL PUNTOS"
$ws.Range("A16:C16").Font.Bold = $true
$ws.Range("A16:C16").HorizontalAlignment = -4108
$ws.Range("A16:C16").VerticalAlignment = -4108
$ws.Range("A16:C16").Borders.LineStyle = 1

# A15:C15 merged "EJEMPLOS"
$ws.Range("A15:C15").Merge()
$ws.Range("A15").Value = "EJEMPLOS"
$ws.Range("A15:C15").Font.Bold = $true
$ws.Range("A15:C15").HorizontalAlignment = -4108
$ws.Range("A15:C15").VerticalAlignment = -4108
$ws.Range("A15:C15").Borders.LineStyle = 1

# Example rows 17-20
$ws.Range("A17").Value = 30
$ws.Range("B17").Value = 8
$ws.Range("C17").Formula = "=A17*B17/5"

$ws.Range("A18").Value = 17.5
$ws.Range("B18").Value = 9
$ws.Range("C18").Formula = "=A18*B18/5"

$ws.Range("A19").Value = 5.5
$ws.Range("B19").Value = 6
$ws.Range("C19").Formula = "=A19*B19/5"

$ws.Range("A20").Value = 5
$ws.Range("B20").Value = 7
$ws.Range("C20").Formula = "=A20*B20/5"

# Page setup (paper size 9 = A4, portrait orientation)
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Selection
$ws.Range("D16").Select()
